$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit Number Info")

# Remove the "pdf_floor_plan_1_7" / "02" entry in row 3 and clear row 4's A/B cells
$ws.Range("A3:B4").ClearContents()
